# Rotate the data held in columns E, F, G (for every row, including the header)
# so that:
#   new E = old G
#   new F = old E
#   new G = old F
#
# This repairs the column ordering for the "group"/"category" metadata
# columns (codeforiati:group-name / codeforiati:category-name /
# codeforiati:group-code -> codeforiati:group-code / codeforiati:group-name /
# codeforiati:category-name) without touching columns A-D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count

$srcRange = $ws.Range("E1:G$rowCount")
$src = $srcRange.Value()

$dst = New-Object 'object[,]' $rowCount,3

for ($i = 1; $i -le $rowCount; $i++) {
    $oldE = $src[$i,1]
    $oldF = $src[$i,2]
    $oldG = $src[$i,3]

    $dst[$i-1,0] = $oldG
    $dst[$i-1,1] = $oldE
    $dst[$i-1,2] = $oldF
}

$ws.Range("E1:G$rowCount").Value = $dst
